# Regenerate merged AHB file:
#  - rename the "_old" / "_new" comparison-column headers to the new
#    version tags "_FV2404" / "_FV2410"
#  - freeze the header row
#  - turn the used range into a native Excel table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. rename header row -------------------------------------------------
# Columns A:J carried the "<name>_old" headers, columns L:U carried the
# matching "<name>_new" headers (column K is the literal "diff" header and
# stays untouched).
$oldHeaderRange = $ws.Range("A1:J1")
[void]$oldHeaderRange.Replace("_old", "_FV2404")

$newHeaderRange = $ws.Range("L1:U1")
[void]$newHeaderRange.Replace("_new", "_FV2410")

# --- 2. freeze the header row ----------------------------------------------
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. convert the used range into a table --------------------------------
$tableRange = $ws.Range("A1:U59")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

Write-Host "done"
